$d = $word.ActiveDocument

$d.Content.Find.Execute("92×84=7728", $true, $false, $false, $false, $false, $true, 1, $false, "27×28=756", 2) | Out-Null
$d.Content.Find.Execute("87×54=4698", $true, $false, $false, $false, $false, $true, 1, $false, "67×76=5092", 2) | Out-Null
$d.Content.Find.Execute("36×53=1908", $true, $false, $false, $false, $false, $true, 1, $false, "98×37=3626", 2) | Out-Null
$d.Content.Find.Execute("32×97=3104", $true, $false, $false, $false, $false, $true, 1, $false, "25×91=2275", 2) | Out-Null
$d.Content.Find.Execute("44×12=528", $true, $false, $false, $false, $false, $true, 1, $false, "32×21=672", 2) | Out-Null
$d.Content.Find.Execute("71×15=1065", $true, $false, $false, $false, $false, $true, 1, $false, "69×51=3519", 2) | Out-Null
$d.Content.Find.Execute("44×24=1056", $true, $false, $false, $false, $false, $true, 1, $false, "97×13=1261", 2) | Out-Null
$d.Content.Find.Execute("96×22=2112", $true, $false, $false, $false, $false, $true, 1, $false, "21×53=1113", 2) | Out-Null
$d.Content.Find.Execute("68×39=2652", $true, $false, $false, $false, $false, $true, 1, $false, "24×21=504", 2) | Out-Null
$d.Content.Find.Execute("16×84=1344", $true, $false, $false, $false, $false, $true, 1, $false, "67×95=6365", 2) | Out-Null
$d.Content.Find.Execute("55×44=2420", $true, $false, $false, $false, $false, $true, 1, $false, "24×93=2232", 2) | Out-Null
$d.Content.Find.Execute("19×72=1368", $true, $false, $false, $false, $false, $true, 1, $false, "90×88=7920", 2) | Out-Null
$d.Content.Find.Execute("90×93=8370", $true, $false, $false, $false, $false, $true, 1, $false, "45×85=3825", 2) | Out-Null
$d.Content.Find.Execute("55×71=3905", $true, $false, $false, $false, $false, $true, 1, $false, "14×75=1050", 2) | Out-Null
$d.Content.Find.Execute("77×65=5005", $true, $false, $false, $false, $false, $true, 1, $false, "39×17=663", 2) | Out-Null
$d.Content.Find.Execute("84×78=6552", $true, $false, $false, $false, $false, $true, 1, $false, "61×17=1037", 2) | Out-Null
$d.Content.Find.Execute("42×79=3318", $true, $false, $false, $false, $false, $true, 1, $false, "13×66=858", 2) | Out-Null
$d.Content.Find.Execute("21×55=1155", $true, $false, $false, $false, $false, $true, 1, $false, "72×37=2664", 2) | Out-Null
$d.Content.Find.Execute("30×27=810", $true, $false, $false, $false, $false, $true, 1, $false, "30×79=2370", 2) | Out-Null
$d.Content.Find.Execute("51×66=3366", $true, $false, $false, $false, $false, $true, 1, $false, "52×98=5096", 2) | Out-Null
$d.Content.Find.Execute("54×15=810", $true, $false, $false, $false, $false, $true, 1, $false, "23×39=897", 2) | Out-Null
$d.Content.Find.Execute("35×34=1190", $true, $false, $false, $false, $false, $true, 1, $false, "27×66=1782", 2) | Out-Null
$d.Content.Find.Execute("52×15=780", $true, $false, $false, $false, $false, $true, 1, $false, "78×63=4914", 2) | Out-Null
$d.Content.Find.Execute("42×33=1386", $true, $false, $false, $false, $false, $true, 1, $false, "86×32=2752", 2) | Out-Null
$d.Content.Find.Execute("86×36=3096", $true, $false, $false, $false, $false, $true, 1, $false, "70×55=3850", 2) | Out-Null

Write-Host "Done replacing multiplication problems."
